# RTO work order list
# Adds "RTO work order" columns (Order Type, Order_Status, Zip, Country,
# City, State) with one sample data row to the RTO_DispatchProcess sheet.

$wb = $excel.ActiveWorkbook

# The sheet is already the active one in the workbook (tabSelected="1"),
# but look it up by name too, to be resilient either way.
$ws = $wb.Worksheets.Item("RTO_DispatchProcess")
$ws.Activate()

# New header row (row 1), columns H:M, entered column-by-column
# (header cell immediately followed by its sample value below it),
# matching how the shared-string table grows in the authored workbook.
$ws.Range("H1").Value = "Order Type"
$ws.Range("H2").Value = "Cleaning"

$ws.Range("I1").Value = "Order_Status"
$ws.Range("I2").Value = "Open"

$ws.Range("J1").Value = "Zip"
$ws.Range("J2").Value = 75024

$ws.Range("K1").Value = "Country"
$ws.Range("K2").Value = "United States"

$ws.Range("L1").Value = "City"
$ws.Range("L2").Value = "Plano"

$ws.Range("M1").Value = "State"
$ws.Range("M2").Value = "TX"

# Autofit the newly populated columns, as Excel does after typing data.
$ws.Columns("G:I").AutoFit()

# Leave the selection where the author's last click landed.
$ws.Range("I14").Select()
